$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")

$ws.Cells.Item(2,4).Value = 2417.44
$ws.Cells.Item(2,5).Value = 100.05
$ws.Cells.Item(3,4).Value = 1925
$ws.Cells.Item(3,5).Value = 630
$ws.Cells.Item(4,4).Value = 1915.88
$ws.Cells.Item(4,5).Value = 633.6
$ws.Cells.Item(6,4).Value = 1715
$ws.Cells.Item(6,5).Value = 550
$ws.Cells.Item(7,4).Value = 1670
$ws.Cells.Item(7,5).Value = 530
$ws.Cells.Item(8,4).Value = 1625
$ws.Cells.Item(9,4).Value = 1090.71
$ws.Cells.Item(9,5).Value = 358.18
$ws.Cells.Item(10,4).Value = 1054.97
$ws.Cells.Item(10,5).Value = 350.03
$ws.Cells.Item(11,1).Value = "BRVM - AGRICULTURE"
$ws.Cells.Item(11,3).Value = 3
$ws.Cells.Item(11,4).Value = 962.36
$ws.Cells.Item(11,5).Value = 317.44
$ws.Cells.Item(12,1).Value = "BRVM - INDUSTRIE"
$ws.Cells.Item(12,4).Value = 789.78
$ws.Cells.Item(12,5).Value = 264.18
$ws.Cells.Item(13,1).Value = "BRVM - CONSOMMATION DE BASE"
$ws.Cells.Item(13,4).Value = 651.68
$ws.Cells.Item(13,5).Value = 217.17
$ws.Cells.Item(14,1).Value = "BRVM-PRINCIPAL"
$ws.Cells.Item(14,4).Value = 566.19
$ws.Cells.Item(14,5).Value = 188.31
$ws.Cells.Item(15,1).Value = "BRVM - INDUSTRIELS"
$ws.Cells.Item(15,4).Value = 412.68
$ws.Cells.Item(15,5).Value = 136.99
$ws.Cells.Item(16,1).Value = "BRVM-PRESTIGE"
$ws.Cells.Item(16,4).Value = 385.62
$ws.Cells.Item(16,5).Value = 128.77
$ws.Cells.Item(17,1).Value = "BRVM - FINANCES"
$ws.Cells.Item(17,4).Value = 363.41
$ws.Cells.Item(17,5).Value = 121.3
$ws.Cells.Item(18,1).Value = "BRVM - SERVICES FINANCIERS"
$ws.Cells.Item(18,4).Value = 357.15
$ws.Cells.Item(18,5).Value = 119.21
$ws.Cells.Item(19,1).Value = "BRVM - ENERGIE"
$ws.Cells.Item(19,4).Value = 324.85
$ws.Cells.Item(19,5).Value = 107
$ws.Cells.Item(20,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws.Cells.Item(20,4).Value = 313.25
$ws.Cells.Item(20,5).Value = 103.68
$ws.Cells.Item(21,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws.Cells.Item(21,4).Value = 275.97
$ws.Cells.Item(21,5).Value = 92.02
$ws.Cells.Item(22,1).Value = "SAFCA CI (SAFC)"
$ws.Cells.Item(22,2).Value = 3
$ws.Cells.Item(22,3).Value = 0
$ws.Cells.Item(22,4).Value = 21.3
$ws.Cells.Item(22,5).Value = 7.08
$ws.Cells.Item(22,6).Value = "🟢 Achat"
$ws.Cells.Item(22,7).Value = "✅ Renforcer"
$ws.Cells.Item(23,1).Value = "BERNABE CI (BNBC)"
$ws.Cells.Item(23,2).Value = 2
$ws.Cells.Item(23,4).Value = 9.33
$ws.Cells.Item(23,5).Value = 7.14
$ws.Cells.Item(23,6).Value = "🟡 Observer"
$ws.Cells.Item(23,7).Value = "➖ Neutre"
$ws.Cells.Item(24,1).Value = "VIVO ENERGY CI (SHEC)"
$ws.Cells.Item(24,4).Value = 8.44
$ws.Cells.Item(24,5).Value = 5.09
$ws.Cells.Item(25,1).Value = "ORAGROUP TOGO (ORGT)"
$ws.Cells.Item(25,4).Value = 6.23
$ws.Cells.Item(25,5).Value = 3.7
$ws.Cells.Item(26,1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws.Cells.Item(26,4).Value = 5.88
$ws.Cells.Item(26,5).Value = 5.88
$ws.Cells.Item(27,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws.Cells.Item(27,4).Value = 4.07
$ws.Cells.Item(27,5).Value = 4.07
$ws.Cells.Item(28,1).Value = "SICABLE CI (CABC)"
$ws.Cells.Item(28,4).Value = 4
$ws.Cells.Item(28,5).Value = 4
$ws.Cells.Item(29,1).Value = "SITAB CI (STBC)"
$ws.Cells.Item(29,4).Value = 3.06
$ws.Cells.Item(29,5).Value = 3.06
$ws.Cells.Item(30,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws.Cells.Item(30,4).Value = 0.42
$ws.Cells.Item(30,5).Value = -6.25
$ws.Cells.Item(31,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws.Cells.Item(31,4).Value = 0.33
$ws.Cells.Item(31,5).Value = 4.04
$ws.Cells.Item(32,1).Value = "TOTAL"
$ws.Cells.Item(32,2).Value = 0
$ws.Cells.Item(32,3).Value = 3
$ws.Cells.Item(32,4).Value = 0
$ws.Cells.Item(32,5).Value = 0
$ws.Cells.Item(32,7).Value = "➖ Neutre"
$ws.Cells.Item(33,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws.Cells.Item(33,3).Value = 1
$ws.Cells.Item(33,4).Value = -1.06
$ws.Cells.Item(33,5).Value = -1.06
$ws.Cells.Item(34,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws.Cells.Item(34,4).Value = -1.32
$ws.Cells.Item(34,5).Value = -1.32
$ws.Cells.Item(35,1).Value = "AIR LIQUIDE CI (SIVC)"
$ws.Cells.Item(35,4).Value = -1.82
$ws.Cells.Item(35,5).Value = -1.82
$ws.Cells.Item(36,1).Value = "FILTISAC CI (FTSC)"
$ws.Cells.Item(36,4).Value = -1.84
$ws.Cells.Item(36,5).Value = -1.84
$ws.Cells.Item(37,1).Value = "CIE CI (CIEC)"
$ws.Cells.Item(37,2).Value = 0
$ws.Cells.Item(37,4).Value = -2.27
$ws.Cells.Item(37,5).Value = -2.27
$ws.Cells.Item(37,7).Value = "➖ Neutre"
$ws.Cells.Item(38,1).Value = "BANK OF AFRICA NG (BOAN)"
$ws.Cells.Item(38,4).Value = -2.71
$ws.Cells.Item(38,5).Value = -2.71
$ws.Cells.Item(39,1).Value = "SAPH CI (SPHC)"
$ws.Cells.Item(39,4).Value = -2.76
$ws.Cells.Item(39,5).Value = -2.76
$ws.Cells.Item(40,1).Value = "SOGB CI (SOGC)"
$ws.Cells.Item(40,4).Value = -3.23
$ws.Cells.Item(40,5).Value = -3.23
$ws.Cells.Item(41,1).Value = "SOLIBRA CI (SLBC)"
$ws.Cells.Item(41,4).Value = -3.23
$ws.Cells.Item(41,5).Value = -3.23
$ws.Cells.Item(42,1).Value = "SETAO CI (STAC)"
$ws.Cells.Item(42,3).Value = 1
$ws.Cells.Item(42,4).Value = -3.45
$ws.Cells.Item(42,5).Value = -3.45
$ws.Cells.Item(44,1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws.Cells.Item(44,2).Value = 0
$ws.Cells.Item(44,3).Value = 1
$ws.Cells.Item(44,4).Value = -6.55
$ws.Cells.Item(44,5).Value = -6.55
$ws.Cells.Item(44,6).Value = "🟡 Observer"
$ws.Cells.Item(44,7).Value = "➖ Neutre"

$wst = $wb.Worksheets.Item("Top_YTD")
$wst.Cells.Item(2,2).Value = 422480.94
$wst.Cells.Item(3,2).Value = 40688.75
$wst.Cells.Item(4,2).Value = 40195.79
$wst.Cells.Item(6,2).Value = 30177
$wst.Cells.Item(7,2).Value = 28174.4
$wst.Cells.Item(8,2).Value = 26319.2
$wst.Cells.Item(9,2).Value = 9860.84
$wst.Cells.Item(10,2).Value = 9113.42
$wst.Cells.Item(11,2).Value = 7350.11
